# Weekly update: two new daily-price rows are published for
# "Fruta, Terminal Hortofrutícola Agro Chillán - Ciruela".
# They get inserted right after row 28 (i.e. become the new rows 29-30),
# which pushes all the existing data rows (old 29-108) down by two rows
# (new 31-110). The sheet's used range therefore grows from A1:T108 to
# A1:T110.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data (rows 29:108) down by two rows, inserting two
# blank rows at 29:30 for the new entries.
$ws.Rows("29:30").Insert()

# --- New row 29 -----------------------------------------------------
$ws.Range("A29").Value = 7
$ws.Range("B29").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C29").Value = "Ñuble"
$ws.Range("D29").Value = 45014
$ws.Range("E29").Value = 16
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100103
$ws.Range("H29").Value = "Frutos de hueso (carozo)"
$ws.Range("I29").Value = 100103002
$ws.Range("J29").Value = "Ciruela"
$ws.Range("K29").Value = "Angeleno"
$ws.Range("L29").Value = "Especial"
$ws.Range("M29").Value = 60
$ws.Range("N29").Value = 12000
$ws.Range("O29").Value = 12000
$ws.Range("P29").Value = 12000
$ws.Range("Q29").Value = "$/bandeja 18 kilos granel"
$ws.Range("R29").Value = "Región de O'Higgins"
$ws.Range("S29").Value = 667
$ws.Range("T29").Value = 18

# --- New row 30 -----------------------------------------------------
$ws.Range("A30").Value = 7
$ws.Range("B30").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C30").Value = "Ñuble"
$ws.Range("D30").Value = 45014
$ws.Range("E30").Value = 16
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100103
$ws.Range("H30").Value = "Frutos de hueso (carozo)"
$ws.Range("I30").Value = 100103002
$ws.Range("J30").Value = "Ciruela"
$ws.Range("K30").Value = "Angeleno"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 50
$ws.Range("N30").Value = 11000
$ws.Range("O30").Value = 11000
$ws.Range("P30").Value = 11000
$ws.Range("Q30").Value = "$/bandeja 18 kilos granel"
$ws.Range("R30").Value = "Región de O'Higgins"
$ws.Range("S30").Value = 611
$ws.Range("T30").Value = 18

# Give the new date cells the same date/time number format as the rest
# of column D.
$ws.Range("D29:D30").NumberFormat = $ws.Range("D31").NumberFormat()
